# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker rows (B16:G24) get re-sorted: the "new" workers block
# (OSCAR / EDEL / YOHON / OSVALDO) moves to the top (rows 16-20) and the
# existing "YOSIMAR LORDUY CHAVEZ" rows move below it (rows 21-24), now
# sorted by ascending period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Target data for rows 16-24 : Tipo(B) Doc(C) Nombre(D) Periodo(E) Mora(F) Salario(G)
$data = @(
    @("CC", "73188736",   "OSCAR EUGENIO ALVAREZ ANGULO",     "1608", 100454, 781242),
    @("CC", "73270962",   "EDEL ENRIQUE CHAMORRO CANTILLO",   "1609", 19200,  781242),
    @("CC", "1049482102", "YOHON DAIRO DE LEON BERROCAL",     "1609", 27734,  781242),
    @("CC", "73270962",   "EDEL ENRIQUE CHAMORRO CANTILLO",   "1610", 32000,  781242),
    @("CC", "9294622",    "OSVALDO ENRIQUE ALVAREZ MARTINEZ", "1701", 96000,  781242),
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ",            "1703", 29509,  737717),
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ",            "1704", 29509,  737717),
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ",            "1705", 29509,  737717),
    @("CC", "1128058778", "YOSIMAR LORDUY CHAVEZ",            "1706", 29509,  737717)
)

$row = 16
foreach ($r in $data) {
    $ws.Range("B$row").Value = $r[0]
    $ws.Range("C$row").Value = $r[1]
    $ws.Range("D$row").Value = $r[2]
    $ws.Range("E$row").Value = $r[3]
    $ws.Range("F$row").Value = $r[4]
    $ws.Range("G$row").Value = $r[5]
    $row++
}
